$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 6
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 5

$ws.Range("H3").Select()
